$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.856.59"
$ws.Range("E2").Value = "  +2.83%  "
$ws.Range("D3").Value = "1.901.24"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.038"
$ws.Range("E4").Value = "  +3.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.74"
$ws.Range("E5").Value = "  +2.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.033"
$ws.Range("E6").Value = "  +3.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5207"
$ws.Range("E7").Value = "  +1.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4022"
$ws.Range("E8").Value = "  +3.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08481"
$ws.Range("E9").Value = "  +1.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.129"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.70"
$ws.Range("E11").Value = "  +3.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.349"
$ws.Range("E12").Value = "  +2.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.82"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("B14").Value = "BinanceUSD"
$ws.Range("C14").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.041"
$ws.Range("E14").Value = "  +3.81%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.333"
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("D16").Value = "1.784.59"
$ws.Range("E16").Value = "  -5.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001125"
$ws.Range("E17").Value = "  +2.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.34"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06878"
$ws.Range("E19").Value = "  +3.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.99"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.033"
$ws.Range("E21").Value = "  +3.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.061"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").Value = "28.877.89"
$ws.Range("E23").Value = "  +2.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.328"
$ws.Range("E25").Value = "  +3.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.91"
$ws.Range("E26").Value = "  +3.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.13"
$ws.Range("E27").Value = "  +2.89%  "
$ws.Range("D28").Value = "1.988.10"
$ws.Range("E28").Value = "  -4.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.417"
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.28"
$ws.Range("E30").Value = "  +3.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1068"
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.059"
$ws.Range("E32").Value = "  +2.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.916"
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.691"
$ws.Range("E34").Value = "  +2.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02473"
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06612"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2217"
$ws.Range("E37").Value = "  +1.75%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.124"
$ws.Range("E38").Value = "  -4.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.282"
$ws.Range("E39").Value = "  +5.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.209"
$ws.Range("E40").Value = "  +1.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6545"
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.089"
$ws.Range("E42").Value = "  +2.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.42"
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6121"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.23"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.784"
$ws.Range("E46").Value = "  +3.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.240"
$ws.Range("E47").Value = "  -3.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.033"
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.230"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "123.04"
$ws.Range("E50").Value = "  +2.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06919"
$ws.Range("E51").Value = "  +0.42%  "
